$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that are missing the "Processed" marker in column C (every 10th row
# starting at 43, through 1033).
for ($r = 43; $r -le 1033; $r += 10) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
